# tests/test-data.xlsx — add SKU header column, refresh quantities, append
# the new 70123_410 size run.  (update readme, fix test bug)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column A (column B already reads "Quantity Available").
$ws.Range("A1").Value = "SKU"

# Refreshed quantities for the existing 70030_200 rows.
$ws.Range("B2").Value = 14
$ws.Range("B4").Value = 93
$ws.Range("B5").Value = 183
$ws.Range("B6").Value = 220
$ws.Range("B7").Value = 204

# Row 10 used to hold the one-off "70115_990-L" SKU; it is replaced by the
# first row of the new 70123_410 run.
$ws.Range("A10").Value = "70123_410-XS"
$ws.Range("B10").Value = 143

# Append the rest of the new 70123_410 size run.
$ws.Range("A11").Value = "70123_410-S"
$ws.Range("B11").Value = 136
$ws.Range("A12").Value = "70123_410-M"
$ws.Range("B12").Value = 300
$ws.Range("A13").Value = "70123_410-L"
$ws.Range("B13").Value = 741
$ws.Range("A14").Value = "70123_410-XL"
$ws.Range("B14").Value = 544
$ws.Range("A15").Value = "70123_410-2XL"
$ws.Range("B15").Value = 235
$ws.Range("A16").Value = "70123_410-3XL"
$ws.Range("B16").Value = 73

# Carry the existing row style down onto the newly-added rows.
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A11:B16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the updated cell-style bookkeeping from the workbook edit.
$builtin = $wb.Styles.Item("Excel Built-in Explanatory Text")
$builtin.Name = "Normal"

# View changes: zoom out and re-select the now-larger data range.
$excel.ActiveWindow.Zoom = 230
$ws.Range("A1:B16").Select() | Out-Null
